$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '36.646.15'
$ws.Range('E2').Value = '  -0.83%  '
$ws.Range('D3').Value = '2.058.29'
$ws.Range('E3').Value = '  +1.00%  '
$ws.Range('E4').Value = '  +0.28%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '243.48'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.60%  '
$ws.Range('E6').Value = '  +1.83%  '
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '54.31'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -6.53%  '
$ws.Range('E9').Value = '  +0.15%  '
$ws.Range('E10').Value = '  -3.08%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0749'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.48%  '
$ws.Range('E12').Value = '  -3.18%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.932'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +6.08%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '14.72'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -3.64%  '
$ws.Range('D15').Value = '2.360.66'
$ws.Range('E15').Value = '  +1.33%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.46'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -2.69%  '
$ws.Range('D17').Value = '2.095.00'
$ws.Range('E17').Value = '  +3.17%  '
$ws.Range('D18').Value = '36.569.80'
$ws.Range('E18').Value = '  -0.84%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '16.97'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -6.69%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '71.98'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.97%  '
$ws.Range('D21').Value = '0.0₃0862'
$ws.Range('E21').Value = '  -2.17%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '238.01'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.20%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.25'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.64%  '
$ws.Range('E24').Value = '  -0.07%  '
$ws.Range('E25').Value = '  -3.19%  '
$ws.Range('E26').Value = '  -0.60%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.33'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.30%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '164.22'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -2.05%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '20.12'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.30%  '
$ws.Range('E30').Value = '  -1.21%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '5.08'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -8.43%  '
$ws.Range('E32').Value = '  +7.78%  '
$ws.Range('E33').Value = '  -4.29%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0596'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -2.12%  '
$ws.Range('E35').Value = '  +0.28%  '
$ws.Range('E36').Value = '  -0.47%  '
$ws.Range('E37').Value = '  -1.46%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0824'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -4.11%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.25'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -3.52%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '4.89'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -5.67%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.89'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -6.72%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0215'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.62%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.10'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.82%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '94.24'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.76%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0909'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -4.07%  '
$ws.Range('D46').Value = '1.408.55'
$ws.Range('E46').Value = '  +9.18%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '7.60'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +14.14%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '15.98'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -5.23%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.91'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.98%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.27'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -3.34%  '
$ws.Range('D51').Value = '2.247.86'
$ws.Range('E51').Value = '  +1.48%  '
